$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the timestamp precision on the existing tail rows (38-45) ---
$ws.Range("A38").Value = 45913.91012901621
$ws.Range("A39").Value = 45913.91039144676
$ws.Range("A40").Value = 45913.91062084491
$ws.Range("A41").Value = 45913.91133267361
$ws.Range("A42").Value = 45913.91175178241
$ws.Range("A43").Value = 45913.91236648148
$ws.Range("A44").Value = 45913.91350438658
$ws.Range("A45").Value = 45913.91387811342

# --- Append new rows 46-48 ---
$ws.Range("A46").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A46").Value = 45913.91664212963
$ws.Range("B46").Value = 660
$ws.Range("C46").Value = 264
$ws.Range("D46").Value = "0:0"
$ws.Range("E46").Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"

$ws.Range("A47").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A47").Value = 45913.92082518519
$ws.Range("B47").Value = 660
$ws.Range("C47").Value = 264
$ws.Range("D47").Value = "0:0"
$ws.Range("E47").Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"

$ws.Range("A48").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A48").Value = 45913.921534596
$ws.Range("B48").Value = 660
$ws.Range("C48").Value = 264
$ws.Range("D48").Value = "0:0"
$ws.Range("E48").Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"
